$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2021" year column to the table, mirroring the formatting of
# the existing last column (Q) for the header row and of the "0.0"-formatted
# data cells in the data row.
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)
$ws.Range("R4").Value = 2021

$ws.Range("H5").Copy()
$ws.Range("R5").PasteSpecial(-4122)
$ws.Range("R5").Value = 18.953297329007047

# Move the active selection (author re-saved with the cursor one row up)
$ws.Range("Q8").Select()
